# Update town close year columns (#769)
# Rename the "2023/2024"-style headers to generic "Prior Year / Curr. Year"
# labels and widen the affected columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text (row 1) ------------------------------------------------
# Columns A-H, Q, R keep their existing headers; only I1:P1 change from
# explicit-year labels to "Prior Year"/"Curr. Year" labels.
$ws.Range("I1").Value = "Curr. Year Card BMV"
$ws.Range("J1").Value = "Prior Year Card BMV"
$ws.Range("K1").Value = "Curr. Year LMV"
$ws.Range("L1").Value = "Curr. Year BMV"
$ws.Range("M1").Value = "Curr. Year Total MV"
$ws.Range("N1").Value = "Prior Year LMV"
$ws.Range("O1").Value = "Prior Year BMV"
$ws.Range("P1").Value = "Prior Year Total MV"

# --- Column widths --------------------------------------------------------
# Columns I:P (9-16) are widened (and lose their "best fit" auto flag,
# since they are now explicitly sized) to accommodate the longer text.
$ws.Columns.Item(9).ColumnWidth = 20.833333333333332
$ws.Columns.Item(10).ColumnWidth = 20.833333333333332
$ws.Columns.Item(11).ColumnWidth = 15.666666666666666
$ws.Columns.Item(12).ColumnWidth = 15.5
$ws.Columns.Item(13).ColumnWidth = 19.333333333333332
$ws.Columns.Item(14).ColumnWidth = 16.0
$ws.Columns.Item(15).ColumnWidth = 16.0
$ws.Columns.Item(16).ColumnWidth = 19.666666666666668

# Keep the default selection/position at A1.
$ws.Range("A1").Select()
